# Applies the "Added flowchart PNGs to assets folder" update to the
# Approvals report: inserts/re-orders several approval rows (new users,
# image filename lists and dates), extending the table from 14 to 20
# data rows (A2:C21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Approvals")

$data = @(
    @('prajyotkankal9@gmail.com', '1748856203722-414725785.png, 1749036756302-383338020.png, 1749036801330-106037180.png', '6/17/2025'),
    @('admin', '1748799431337-109774176.jpg, 1748798941789-920707644.jpg', '6/13/2025'),
    @('prajyotkankal9@gmail.com', '1748799431337-109774176.jpg, 1748798941789-920707644.jpg', '6/12/2025'),
    @('admin', '1749036756302-383338020.png', '6/10/2025'),
    @('admin', '1748798941789-920707644.jpg', '6/10/2025'),
    @('admin', '1749030981363-555614522.png, 1749036756302-383338020.png', '6/10/2025'),
    @('admin', '1749036756302-383338020.png, 1749036801330-106037180.png', '6/10/2025'),
    @('prajyotkankal9@gmail.com', '1749030981363-555614522.png, 1749036740156-172730154.png', '6/10/2025'),
    @('admin', '1748798941789-920707644.jpg', '6/9/2025'),
    @('prajyotkankal9@gmail.com', '1748798941789-920707644.jpg', '6/9/2025'),
    @('prajyotkankal9@gmail.com', '1748798941789-920707644.jpg', '6/9/2025'),
    @('sanketnk1401@gmail.com', '1749036723774-181764139.png, 1749036740156-172730154.png, 1749036756302-383338020.png', '6/7/2025'),
    @('prajyotkankal9@gmail.com', '1749030981363-555614522.png, 1749036740156-172730154.png, 1749036756302-383338020.png', '6/6/2025'),
    @('prajyotkankal9@gmail.com', '1748856203722-414725785.png, 1749030981363-555614522.png', '6/5/2025'),
    @('prajyotkankal9@gmail.com', '1749036723774-181764139.png, 1749036740156-172730154.png, 1749036756302-383338020.png', '6/5/2025'),
    @('prajyotkankal9@gmail.com', '1748798941789-920707644.jpg, 1748799431337-109774176.jpg', '6/4/2025'),
    @('prajyotkankal12@gmail.com', '1748798941789-920707644.jpg', '6/2/2025'),
    @('prajyotkankal9@gmail.com', '1748799431337-109774176.jpg, 1748798941789-920707644.jpg', '6/2/2025'),
    @('prajyotkankal9@gmail.com', '1748798941789-920707644.jpg, 1748799431337-109774176.jpg', '6/2/2025'),
    @('prajyotkankal9@gmail.com', '1748856203722-414725785.png', '6/2/2025'),
)

# First pass: write every cell so the sheet grows to A1:C21. Column C
# values look like dates, so Excel will auto-convert them to date serial
# numbers with a date number format here; that gets corrected below.
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# Second pass: force column C back to plain text matching the source
# report (e.g. "6/17/2025" as a literal string, not a formatted date),
# then drop the temporary text number format so cells keep the default
# (unstyled) cell style, same as the rest of the sheet.
$dateRange = $ws.Range("C2:C21")
$dateRange.NumberFormat = "@"
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
$dateRange.ClearFormats()
